# Update cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.065.16"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.16%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.384.20"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "574.33"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "137.45"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.381.49"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.20%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -1.16%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.65"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.95%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.122"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.34%  "
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.64%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.960.48"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.21%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +0.77%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -3.21%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "25.76"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.380.99"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.82%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.194.12"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.85"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.75"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -1.22%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.31"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -1.83%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "377.01"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.35%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -2.38%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.517.84"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.49%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000126"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -1.42%  "
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.70%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +12.34%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.67"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.39%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.41"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -2.80%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.83%  "
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.64%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.43"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.15"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -4.41%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.55"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.73%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.83"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -0.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "164.51"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0758"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -3.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "25.79"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +2.84%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.775"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -1.44%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.13%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.02%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -4.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.543.31"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +8.32%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.80"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -1.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "22.93"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.44"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +3.79%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0259"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -1.47%  "
